# Recodificación de funciones (formato):
#  - Tienda "Velázquez" (con TPVs "BAR" / "SERVIDOR TIENDA") se unifica en
#    un único punto de venta "LOCAL LM" (columna Nombre_TPV desaparece).
#  - Los importes/operaciones quedan redistribuidos por turno/medio de
#    pago, con dos filas nuevas al final.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renombrar la hoja y la tienda
$ws.Name = "LOCAL LM"

$fecha = $ws.Range("C2").Value2

$rows = @(
    @{ R = 2;  D = "Mañana"; E = "GLOVO";         F = 190.6;   G = 6 },
    @{ R = 3;  D = "Mañana"; E = "SMS";            F = 7.2;     G = 2 },
    @{ R = 4;  D = "Mañana"; E = "EUROS";          F = 1915.04; G = 182 },
    @{ R = 5;  D = "Mañana"; E = "TARJETA VISA";   F = 4663.5;  G = 414 },
    @{ R = 6;  D = "Tarde";  E = "EUROS";          F = 295.5;   G = 32 },
    @{ R = 7;  D = "Tarde";  E = "TARJETA VISA";   F = 857.99;  G = 71 },
    @{ R = 8;  D = "Mañana"; E = "EUROS";          F = 563.28;  G = 46 },
    @{ R = 9;  D = "Mañana"; E = "GLOVO";          F = 73.8;    G = 5 },
    @{ R = 10; D = "Mañana"; E = "TARJETA VISA";   F = 1351.34; G = 90 },
    @{ R = 11; D = "Tarde";  E = "TARJETA VISA";   F = 2196.68; G = 198 },
    @{ R = 12; D = "Tarde";  E = "EUROS";          F = 1107.08; G = 100 },
    @{ R = 13; D = "Tarde";  E = "EUROS";          F = 563.28;  G = 46 },
    @{ R = 14; D = "Tarde";  E = "GLOVO";          F = 73.8;    G = 5 },
    @{ R = 15; D = "Tarde";  E = "TARJETA VISA";   F = 1351.34; G = 90 },
    @{ R = 16; D = "Mañana"; E = "EUROS";          F = 295.5;   G = 32 },
    @{ R = 17; D = "Mañana"; E = "TARJETA VISA";   F = 857.99;  G = 71 }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = "LOCAL LM"

    # Columna B (Nombre_TPV) desaparece para todas las filas de datos
    $ws.Cells.Item($r, 2).ClearContents()

    $cC = $ws.Cells.Item($r, 3)
    $cC.NumberFormat = "DD/MM/YYYY"
    $cC.Value = $fecha

    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    $cF = $ws.Cells.Item($r, 6)
    $cF.NumberFormat = "#,##0.00"
    $cF.Value = $row.F

    $cG = $ws.Cells.Item($r, 7)
    $cG.NumberFormat = "#,##0"
    $cG.Value = $row.G
}
